# Update the "Metadata" sheet: URL, Version, Date and Publisher fields.
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-report-evidence-text"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Clear the stale ele-1/ext-1 constraint text that had leaked onto the
# root "Extension" row's Constraint(s) column in the "Elements" sheet.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").ClearContents()
